$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chests")

# Fill in the new chest rows for "Your Cave" upper level.
# Row order chosen so the shared-string table is appended in the
# same sequence as the target workbook: row 134, 135, 133, 136, 137.

$ws.Range("A4").Value = 134
$ws.Range("B4").Value = "Your Cave (459)"
$ws.Range("C4").Value = "1x Levitation, 1x Healing Potion II"

$ws.Range("A5").Value = 135
$ws.Range("B5").Value = "Your Cave (460)"
$ws.Range("C5").Value = "3x Torch"

$ws.Range("A3").Value = 133
$ws.Range("B3").Value = "Your Cave (459)"
$ws.Range("C3").Value = "2x Rope"

$ws.Range("A6").Value = 136
$ws.Range("B6").Value = "Your Cave (459)"
$ws.Range("C6").Value = "1x Rope"

$ws.Range("A7").Value = 137
$ws.Range("B7").Value = "Your Cave (459)"
$ws.Range("C7").Value = "1x Rope"

# Select C8 (next empty row) as the last active cell, matching the author's
# final cursor position, and make the Chests sheet the active tab.
$ws.Range("C8").Select()
$ws.Activate()
